# Update "想去人数" (number of people wanting to go) figures that were
# refreshed by the data-generation job (commit: "Update gh-pages to output
# generated at 456a3b4").
#
# Sheet "展览" and sheet "全部类型" contain the same rows (the latter is an
# aggregate view across all event types), so both need the identical update:
#   F2: 1702 -> 1703
#   F3: 7851 -> 7855
#   F5: 242  -> 243

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1703
    $ws.Range("F3").Value = 7855
    $ws.Range("F5").Value = 243
}
